# The sheets each had a leading helper/index column (STDID row-number
# column A, styled like the header) that was only used for bookkeeping.
# Having to keep that index in sync made it impossible to just append a
# new row to the "random list" table. Removing the column lets every
# other column shift left (B->A, C->B, ...) and the data start cleanly
# at column A again.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("A1").EntireColumn.Delete()
}
